$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 231.77777
$ws.Range("I5").Value = 82.833336
$ws.Range("J5").Value = 529.6667
$ws.Range("K5").Value = 82.833336
$ws.Range("L5").Value = 529.6667
$ws.Range("M5").Value = 32.166664
$ws.Range("N5").Value = -759.6667
$ws.Range("H15").Value = 300.05405
$ws.Range("I15").Value = 300.05405
$ws.Range("K15").Value = 900.1621500000001
$ws.Range("M15").Value = -731.1621500000001
$ws.Range("H39").Value = 889.1429000000001
$ws.Range("I39").Value = 258.5
$ws.Range("J39").Value = 1730
$ws.Range("K39").Value = 775.5
$ws.Range("L39").Value = 5190
$ws.Range("M39").Value = -479.5
$ws.Range("N39").Value = -5782
$ws.Range("H70").Value = 2640
$ws.Range("I70").Value = 2350
$ws.Range("K70").Value = 7050
$ws.Range("M70").Value = -6780
$ws.Range("H73").Value = 2640
$ws.Range("I73").Value = 2350
$ws.Range("K73").Value = 7050
$ws.Range("M73").Value = -6114
$ws.Range("H98").Value = 894.94446
$ws.Range("I98").Value = 1001.4545
$ws.Range("J98").Value = 727.5714
$ws.Range("K98").Value = 1001.4545
$ws.Range("L98").Value = 727.5714
$ws.Range("M98").Value = 496.5454999999999
$ws.Range("N98").Value = -3723.5714
$ws.Range("H122").Value = 894.94446
$ws.Range("I122").Value = 1001.4545
$ws.Range("J122").Value = 727.5714
$ws.Range("K122").Value = 3004.3635
$ws.Range("L122").Value = 2182.7142
$ws.Range("M122").Value = -554.3635000000004
$ws.Range("N122").Value = -7082.7142
$ws.Range("H137").Value = 2033.5714
$ws.Range("I137").Value = 1429.1364
$ws.Range("J137").Value = 4249.8335
$ws.Range("K137").Value = 4287.4092
$ws.Range("L137").Value = 12749.5005
$ws.Range("M137").Value = -1737.4092
$ws.Range("N137").Value = -17849.5005
$ws.Range("H141").Value = 2279
$ws.Range("I141").Value = 2421.111
$ws.Range("J141").Value = 1000
$ws.Range("K141").Value = 7263.333
$ws.Range("L141").Value = 3000
$ws.Range("M141").Value = -2083.333
$ws.Range("N141").Value = -13360

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 16744.5
$ws.Range("J27").Value = 16744.5
$ws.Range("L27").Value = 16744.5
$ws.Range("N27").Value = -17112.5
$ws.Range("H32").Value = 6279.6665
$ws.Range("I32").Value = 4548.5557
$ws.Range("K32").Value = 4548.5557
$ws.Range("M32").Value = -4261.5557
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 812.5
$ws.Range("I24").Value = 812.5
$ws.Range("K24").Value = 812.5
$ws.Range("M24").Value = -577.5
$ws.Range("H29").Value = 2408
$ws.Range("J29").Value = 1800
$ws.Range("L29").Value = 1800
$ws.Range("N29").Value = -2378
$ws.Range("H36").Value = 7150
$ws.Range("I36").Value = 2980
$ws.Range("J36").Value = 28000
$ws.Range("K36").Value = 2980
$ws.Range("L36").Value = 28000
$ws.Range("M36").Value = -2446
$ws.Range("N36").Value = -29068
$ws.Range("H86").Value = 6931.7
$ws.Range("I86").Value = 5766.7144
$ws.Range("J86").Value = 9650
$ws.Range("K86").Value = 5766.7144
$ws.Range("L86").Value = 9650
$ws.Range("M86").Value = -4643.7144
$ws.Range("N86").Value = -11896
$ws.Range("H89").Value = 6931.7
$ws.Range("I89").Value = 5766.7144
$ws.Range("J89").Value = 9650
$ws.Range("K89").Value = 28833.572
$ws.Range("L89").Value = 48250
$ws.Range("M89").Value = -23217.572
$ws.Range("N89").Value = -59482
$ws.Range("H105").Value = 2058.4285
$ws.Range("I105").Value = 1443.1666
$ws.Range("K105").Value = 1443.1666
$ws.Range("M105").Value = 303.8334
$ws.Range("H134").Value = 4687.4287
$ws.Range("I134").Value = 4968.6665
$ws.Range("K134").Value = 14905.9995
$ws.Range("M134").Value = -12370.9995

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2709.1428
$ws.Range("I31").Value = 2729
$ws.Range("J31").Value = 2590
$ws.Range("K31").Value = 2729
$ws.Range("L31").Value = 2590
$ws.Range("M31").Value = -2434
$ws.Range("N31").Value = -3180
$ws.Range("H34").Value = 2709.1428
$ws.Range("I34").Value = 2729
$ws.Range("J34").Value = 2590
$ws.Range("K34").Value = 2729
$ws.Range("L34").Value = 2590
$ws.Range("M34").Value = -2527
$ws.Range("N34").Value = -2994
$ws.Range("H35").Value = 23500
$ws.Range("I35").Value = 20000
$ws.Range("J35").Value = 27000
$ws.Range("K35").Value = 20000
$ws.Range("L35").Value = 27000
$ws.Range("M35").Value = -19706
$ws.Range("N35").Value = -27588
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H122").Value = 1713
$ws.Range("I122").Value = 998.5
$ws.Range("J122").Value = 2232.6365
$ws.Range("K122").Value = 2995.5
$ws.Range("L122").Value = 6697.9095
$ws.Range("M122").Value = -545.5
$ws.Range("N122").Value = -11597.9095
$ws.Range("H132").Value = 1794.7931
$ws.Range("I132").Value = 1812.9259
$ws.Range("K132").Value = 5438.7777
$ws.Range("M132").Value = -2908.7777

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 60
$ws.Range("I14").Value = 60
$ws.Range("K14").Value = 180
$ws.Range("M14").Value = -7
$ws.Range("H16").Value = 1
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H39").Value = 449.5
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H68").Value = 290.7143
$ws.Range("I68").Value = 289
$ws.Range("K68").Value = 867
$ws.Range("M68").Value = -56
$ws.Range("H69").Value = 2769.1538
$ws.Range("J69").Value = 3999.6667
$ws.Range("L69").Value = 11999.0001
$ws.Range("N69").Value = -13621.0001
$ws.Range("H71").Value = 290.7143
$ws.Range("I71").Value = 289
$ws.Range("K71").Value = 2601
$ws.Range("M71").Value = 1455
$ws.Range("H72").Value = 2769.1538
$ws.Range("J72").Value = 3999.6667
$ws.Range("L72").Value = 35997.0003
$ws.Range("N72").Value = -44109.0003
$ws.Range("H117").Value = 730
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 1000
$ws.Range("I36").Value = 1000
$ws.Range("K36").Value = 1000
$ws.Range("M36").Value = -515
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H70").Value = 50001750
$ws.Range("I70").Value = 50001750
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 50001750
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -50001480
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 50001750
$ws.Range("I73").Value = 50001750
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 50001750
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -50000814
$ws.Range("N73").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H102").Value = 1175.4445
$ws.Range("I102").Value = 1035.8
$ws.Range("J102").Value = 1350
$ws.Range("K102").Value = 1035.8
$ws.Range("L102").Value = 1350
$ws.Range("M102").Value = 586.2
$ws.Range("N102").Value = -4594
$ws.Range("H113").Value = 2200
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3833.3333
$ws.Range("I22").Value = 1750
$ws.Range("K22").Value = 1750
$ws.Range("M22").Value = -1455
$ws.Range("H27").Value = 3833.3333
$ws.Range("I27").Value = 1750
$ws.Range("K27").Value = 1750
$ws.Range("M27").Value = -1643
$ws.Range("H31").Value = 4603
$ws.Range("J31").Value = 4000
$ws.Range("L31").Value = 4000
$ws.Range("N31").Value = -4496
$ws.Range("H55").Value = 563.8
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 563.8
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 563.8
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -909.8

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 4999.25
$ws.Range("J55").Value = 6499
$ws.Range("L55").Value = 6499
$ws.Range("N55").Value = -7053
$ws.Range("H107").Value = 289.75
$ws.Range("I107").Value = 328.2857
$ws.Range("K107").Value = 984.8571000000001
$ws.Range("M107").Value = 935.1428999999999

